$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant rows to append after the current last row (176)
$newRows = @(
    @("Hamroyeva Nigora Husan Qizi",
      "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik",
      "AD3548929", "770", "Navoiy viloyati", "Xatirchi tumani",
      "998934319253", "23-11-2024", "+998934319253"),
    @("Saxibova Muxayyo Saidjonovna",
      "Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik",
      "AD7351130", "771", "Toshkent viloyati", "Qibray tumani",
      "998998313914", "23-11-2024", "+998505003914"),
    @("Mahmudova Dilnoza Xolboy qizi",
      "Defektologiya (logopediya) 576 soatlik",
      "AB7598097", "772", "Jizzax viloyati", "Gʻallaorol tumani",
      "998941603726", "26-11-2024", "+998941603726"),
    @("Usarova Tursunoy Umarovna",
      "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik",
      "AA9143291", "773", "Jizzax viloyati", "Sharof Rashidov tumani",
      "998503014091", "26-11-2024", "+998503014091")
)

$startRow = 177
$endRow = $startRow + $newRows.Count - 1

# Force the incoming values to be stored as text (the source data keeps
# values such as "770", "998934319253" and "23-11-2024" as plain text,
# not numbers/dates), then drop the temporary number format again so the
# cells end up with the same (default) style as the rest of the sheet.
$fillRange = $ws.Range("A$startRow`:I$endRow")
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

$fillRange.ClearFormats()
